$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation is inserted at row 11 (right after the
# already-unchanged header + rows 2-10). Every existing data row from the
# old row 11 through the old row 116 shifts down by one row (to rows
# 12-117); the freshly inserted row 11 receives the brand-new record.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new record's data.
$ws.Range("A11").Value = 4
$ws.Range("B11").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C11").Value = "Los Lagos"
$ws.Range("D11").Value = 44473
$ws.Range("E11").Value = 10
$ws.Range("F11").Value = 100112039
$ws.Range("G11").Value = "Ciboulette"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 80
$ws.Range("K11").Value = 3000
$ws.Range("L11").Value = 3000
$ws.Range("M11").Value = 3000
$ws.Range("N11").Value = "$/docena de atados"
$ws.Range("O11").Value = "Región Metropolitana"
$ws.Range("P11").Value = 1000
$ws.Range("Q11").Value = 3
$ws.Range("R11").Value = "Hortaliza"
